$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- First block: simple value replacements in rows 1-12 ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(4).Cells.Item(1).Range.Text = "271"
$t.Rows.Item(6).Cells.Item(1).Range.Text = "0.00266"
$t.Rows.Item(7).Cells.Item(1).Range.Text = "0.00013"
$t.Rows.Item(8).Cells.Item(1).Range.Text = "0.00008"
$t.Rows.Item(9).Cells.Item(1).Range.Text = "0.00022"
$t.Rows.Item(10).Cells.Item(1).Range.Text = "0.00023"
$t.Rows.Item(11).Cells.Item(1).Range.Text = "0.00026"
$t.Rows.Item(12).Cells.Item(1).Range.Text = "0.04128"

# --- Second block: collapse tab-separated multi-run rows into single values ---
$t.Rows.Item(44).Cells.Item(1).Range.Text = "99.94"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.04"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "72"
